# Fruta / hortaliza, semanal
# Insert two new weekly rows (2022-05-25, serial 44706) for
# "Comercializadora del Agro de Limarí" - Ciruela - Angeleno (Primera/Segunda),
# pushing the existing rows 64-70 down to 66-72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 64 (shifts old 64..70 -> 66..72)
$ws.Rows.Item(64).Insert()
$ws.Rows.Item(64).Insert()

# New row 64: Angeleno / Primera
$ws.Range("A64").Value = 2
$ws.Range("B64").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44706
$ws.Range("E64").Value = 4
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100103
$ws.Range("H64").Value = "Frutos de hueso (carozo)"
$ws.Range("I64").Value = 100103002
$ws.Range("J64").Value = "Ciruela"
$ws.Range("K64").Value = "Angeleno"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 16
$ws.Range("N64").Value = 190000
$ws.Range("O64").Value = 200000
$ws.Range("P64").Value = 195000
$ws.Range("Q64").Value = "$/bins (450 kilos)"
$ws.Range("R64").Value = "Región de O'Higgins"
$ws.Range("S64").Value = 433
$ws.Range("T64").Value = 450

# New row 65: Angeleno / Segunda
$ws.Range("A65").Value = 2
$ws.Range("B65").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 44706
$ws.Range("E65").Value = 4
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100103
$ws.Range("H65").Value = "Frutos de hueso (carozo)"
$ws.Range("I65").Value = 100103002
$ws.Range("J65").Value = "Ciruela"
$ws.Range("K65").Value = "Angeleno"
$ws.Range("L65").Value = "Segunda"
$ws.Range("M65").Value = 16
$ws.Range("N65").Value = 150000
$ws.Range("O65").Value = 160000
$ws.Range("P65").Value = 155000
$ws.Range("Q65").Value = "$/bins (450 kilos)"
$ws.Range("R65").Value = "Región de O'Higgins"
$ws.Range("S65").Value = 344
$ws.Range("T65").Value = 450
